# Convert the "NN%" inline-string percentages in columns H (KAST) and J (HS%)
# into real numeric decimal values (e.g. "45%" -> 0.45), matching the rest
# of the sheet's numeric percentage columns (like column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in @("H", "J")) {
        $cell = $ws.Range("$col$row")
        $raw = $cell.Value2

        if ($raw -ne $null -and $raw.ToString().EndsWith('%')) {
            $num = [double]($raw.ToString().TrimEnd('%')) / 100
            $cell.Value = $num
        }
    }
}
